# Modifications to the SDMX modelling
# - REF_AREA (row 4) codelist changes from CL_AREA to CL_COM_GEO_PICT_L123
# - The four fishing-location dimensions (rows 8-11) switch their codelist
#   column from their bespoke CL_FISHING_* codelists to the shared
#   CL_COM_YESNO codelist, and flip the "Coded list?" flag from N to Y
# - Update the selected cell on the DSD sheet to F9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DSD")

# Fishing location dimensions: rows 8 - 11
# Column F -> CL_COM_YESNO, column G -> Y
$ws.Range("F8").Value = "CL_COM_YESNO"
$ws.Range("G8").Value = "Y"

$ws.Range("F9").Value = "CL_COM_YESNO"
$ws.Range("G9").Value = "Y"

$ws.Range("F10").Value = "CL_COM_YESNO"
$ws.Range("G10").Value = "Y"

$ws.Range("F11").Value = "CL_COM_YESNO"
$ws.Range("G11").Value = "Y"

# REF_AREA row: column F codelist
$ws.Range("F4").Value = "CL_COM_GEO_PICT_L123"

# Update the active selection on the DSD sheet
$ws.Activate()
$ws.Range("F9").Select()
